# call_job.xlsx — rename config sheet to configuration, plus the
# accompanying "questions" sheet content/formatting changes that shipped
# in the same commit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Rename the "configures" sheet to "configuration"
# ---------------------------------------------------------------------
$wsConfig = $wb.Worksheets.Item("configures")
$wsConfig.Name = "configuration"

# ---------------------------------------------------------------------
# 2) questions sheet: insert a new reminder question as row 2, wrap
#    column A, and restore the per-row heights for the two long entries
#    (done before the "greeting" edit below so the new shared strings
#    land in the same order as the refreshed workbook)
# ---------------------------------------------------------------------
$wsQuestions = $wb.Worksheets.Item("questions")

$wsQuestions.Rows.Item(2).Insert()
$wsQuestions.Range("A2").Value = "This is a reminder to let you know that the book report assignment is due this Friday, 3 June 2020.  Are you ready?"
$wsQuestions.Range("B2").Value = "Yes/No"
# the row-insert copies row 1's fill/shading onto B2 - put it back to the
# unstyled look every other "B" cell in the column has
$wsQuestions.Range("B2").Style = "Normal"

# Wrap the long question text in column A (keeps each row's own
# font/fill, just adds wrapText like the refreshed workbook does)
$wsQuestions.Columns.Item(1).WrapText = $true

# The two multi-line questions get an explicit row height in the
# refreshed file
$wsQuestions.Rows.Item(2).RowHeight = 30
$wsQuestions.Rows.Item(4).RowHeight = 30

# ---------------------------------------------------------------------
# 3) configuration sheet: shorten the "greeting" template text
# ---------------------------------------------------------------------
$wsConfig.Range("B2").Value = "Hi {{ username }}, "

# ---------------------------------------------------------------------
# 4) Selection/active-sheet bookkeeping to match the refreshed file:
#    "configuration" becomes the active tab, "questions" keeps its
#    selection on A2, and "receivers" is no longer the active tab.
# ---------------------------------------------------------------------
$wsReceivers = $wb.Worksheets.Item("receivers")
$wsReceivers.Range("B3").Select()

$wsQuestions.Range("A2").Select()

$wsConfig.Activate()
$wsConfig.Range("B28").Select()
